# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
# For each changed row we update Price (D) and Volume(1h) (E); two rows
# (40/41 and 50/51) also swap which coin (B/C) occupies that row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.905.76'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").Value = '3.691.17'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '650.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.71'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +1.88%  '

$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.19'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.447'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.61%  '

$ws.Range("E12").Value = '  -0.33%  '

$ws.Range("D13").Value = '4.316.05'
$ws.Range("E13").Value = '  +0.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.89'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").Value = '3.699.02'
$ws.Range("E15").Value = '  +0.44%  '

$ws.Range("D16").Value = '69.904.32'
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.17'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.54'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.53'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +7.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '473.48'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.57%  '

$ws.Range("E22").Value = '  +0.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.98'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.59%  '

$ws.Range("D24").Value = '3.837.59'
$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.01'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.29%  '

$ws.Range("E28").Value = '  +0.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.67'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.46%  '

$ws.Range("E30").Value = '  -0.40%  '

$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.57'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("E33").Value = '  -0.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.94'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("E35").Value = '  +1.60%  '

$ws.Range("D36").Value = '3.688.03'
$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.50'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.37%  '

$ws.Range("E38").Value = '  -0.14%  '

$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '180.26'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +6.95%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.91'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0907'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.932'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '29.41'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.11'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.76'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000269'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.26%  '

$ws.Range("E49").Value = '  -1.49%  '

$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.26'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.25%  '

$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.87'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.16%  '
